$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "46.815.05"
$ws.Range("E2").Value = "  +4.76%  "
$ws.Range("D3").Value = "2.473.42"
$ws.Range("E3").Value = "  +2.24%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "324.01"
$ws.Range("E5").Value = "  +2.52%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "104.90"
$ws.Range("E6").Value = "  +3.64%  "
$ws.Range("E7").Value = "  +1.75%  "
$ws.Range("E8").Value = "  -0.11%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.540"
$ws.Range("E9").Value = "  +2.09%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.10"
$ws.Range("E10").Value = "  +2.16%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0814"
$ws.Range("E11").Value = "  +2.05%  "
$ws.Range("E12").Value = "  +0.79%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "18.22"
$ws.Range("E13").Value = "  -2.69%  "
$ws.Range("E14").Value = "  +3.00%  "
$ws.Range("D15").Value = "2.858.11"
$ws.Range("E15").Value = "  +2.15%  "
$ws.Range("D16").Value = "2.484.90"
$ws.Range("E16").Value = "  +2.53%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.844"
$ws.Range("E17").Value = "  +1.68%  "
$ws.Range("D18").Value = "46.704.24"
$ws.Range("E18").Value = "  +4.92%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.59"
$ws.Range("E19").Value = "  +2.17%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.46"
$ws.Range("E20").Value = "  +1.54%  "
$ws.Range("E21").Value = "  +2.23%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "70.57"
$ws.Range("E22").Value = "  +2.75%  "
$ws.Range("B23").Value = "ImmutableX"
$ws.Range("C23").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.39"
$ws.Range("E23").Value = "  +5.03%  "
$ws.Range("B24").Value = "BitcoinCash"
$ws.Range("C24").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "249.55"
$ws.Range("E24").Value = "  +2.90%  "
$ws.Range("E25").Value = "  +2.61%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.15"
$ws.Range("E26").Value = "  +3.94%  "
$ws.Range("E27").Value = "  +0.08%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.21"
$ws.Range("E28").Value = "  +0.86%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.82"
$ws.Range("E29").Value = "  +3.66%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "35.06"
$ws.Range("E30").Value = "  +4.29%  "
$ws.Range("E31").Value = "  +5.99%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "49.53"
$ws.Range("E32").Value = "  +1.87%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "19.62"
$ws.Range("E33").Value = "  +0.59%  "
$ws.Range("E34").Value = "  +3.03%  "
$ws.Range("E35").Value = "  -0.03%  "
$ws.Range("E36").Value = "  -0.56%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.61"
$ws.Range("E37").Value = "  +3.40%  "
$ws.Range("E38").Value = "  +1.58%  "
$ws.Range("E39").Value = "  +3.73%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "123.12"
$ws.Range("E40").Value = "  +0.13%  "
$ws.Range("E41").Value = "  +1.86%  "
$ws.Range("E42").Value = "  +0.98%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "20.65"
$ws.Range("E43").Value = "  -2.01%  "
$ws.Range("E44").Value = "  +1.18%  "
$ws.Range("D45").Value = "1.974.12"
$ws.Range("E45").Value = "  +1.80%  "
$ws.Range("E46").Value = "  +1.31%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.10"
$ws.Range("E47").Value = "  -0.63%  "
$ws.Range("E48").Value = "  +3.71%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.00"
$ws.Range("E49").Value = "  -2.48%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "5.37"
$ws.Range("E50").Value = "  +16.98%  "
$ws.Range("E51").Value = "  +4.96%  "
